$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 251 (shifts existing rows 251-400 down to 252-401,
# and extends the used range to A1:R401, mirroring Excel's native
# "Insert Row" behaviour incl. carrying the row-above's cell formatting).
$ws.Rows.Item(251).Insert()

# Populate the newly inserted row with the new weekly price record.
$ws.Cells.Item(251, 1).Value = 9
$ws.Cells.Item(251, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(251, 3).Value = "Metropolitana"
$ws.Cells.Item(251, 4).Value = "03/16/2023"
$ws.Cells.Item(251, 5).Value = 13
$ws.Cells.Item(251, 6).Value = 300000001
$ws.Cells.Item(251, 7).Value = "Rabanito"
$ws.Cells.Item(251, 8).Value = "Sin especificar"
$ws.Cells.Item(251, 9).Value = "Primera"
$ws.Cells.Item(251, 10).Value = 7000
$ws.Cells.Item(251, 11).Value = 3000
$ws.Cells.Item(251, 12).Value = 3000
$ws.Cells.Item(251, 13).Value = 3000
$ws.Cells.Item(251, 14).Value = "$/cien unidades (volumen en unidades)"
$ws.Cells.Item(251, 15).Value = "Provincia de Chacabuco"
$ws.Cells.Item(251, 16).Value = 30
$ws.Cells.Item(251, 17).Value = 100
$ws.Cells.Item(251, 18).Value = "Hortaliza"
